$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must stay text (matches source data format).
# Force text number-format first so Excel does not auto-convert them to actual numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply all the updated cell values (prices, volumes, and the two row swaps).
$ws.Range("D2").Value = "68.399.64"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "3.922.59"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "485.02"
$ws.Range("E5").Value = "  +5.34%  "
$ws.Range("D6").Value = "148.13"
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("D7").Value = "0.620"
$ws.Range("E7").Value = "  -0.78%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -3.34%  "
$ws.Range("E10").Value = "  +8.31%  "
$ws.Range("E11").Value = "  +11.90%  "
$ws.Range("D12").Value = "42.47"
$ws.Range("E12").Value = "  -3.23%  "
$ws.Range("D13").Value = "10.54"
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("D14").Value = "4.554.12"
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("B15").Value = "Uniswap"
$ws.Range("C15").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D15").Value = "14.60"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.928.85"
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "19.71"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("D20").Value = "68.561.97"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("D21").Value = "431.34"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").Value = "3.34"
$ws.Range("E22").Value = "  +2.17%  "
$ws.Range("D23").Value = "14.49"
$ws.Range("E23").Value = "  -2.53%  "
$ws.Range("D24").Value = "86.87"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "11.27"
$ws.Range("E25").Value = "  +12.50%  "
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").Value = "10.48"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("D28").Value = "38.09"
$ws.Range("E28").Value = "  +1.55%  "
$ws.Range("D29").Value = "5.89"
$ws.Range("E29").Value = "  +6.99%  "
$ws.Range("D30").Value = "715.45"
$ws.Range("E30").Value = "  -5.11%  "
$ws.Range("D31").Value = "13.23"
$ws.Range("E31").Value = "  -4.05%  "
$ws.Range("E32").Value = "  -4.96%  "
$ws.Range("D33").Value = "2.82"
$ws.Range("E33").Value = "  +3.58%  "
$ws.Range("D34").Value = "0.0₃0895"
$ws.Range("E34").Value = "  +32.55%  "
$ws.Range("D35").Value = "41.56"
$ws.Range("E35").Value = "  -5.06%  "
$ws.Range("D36").Value = "58.57"
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("E37").Value = "  -6.81%  "
$ws.Range("D38").Value = "5.49"
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  +7.59%  "
$ws.Range("E41").Value = "  -1.93%  "
$ws.Range("D42").Value = "3.06"
$ws.Range("E42").Value = "  +11.52%  "
$ws.Range("E43").Value = "  +1.77%  "
$ws.Range("E44").Value = "  -3.80%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("B47").Value = "LidoDAOToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D47").Value = "3.41"
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "2.16"
$ws.Range("E48").Value = "  +1.75%  "
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").Value = "147.07"
$ws.Range("E50").Value = "  +1.76%  "
$ws.Range("D51").Value = "2.83"
